# "continue switch instead of pause switch"
#
# The underlying commit only touches PowerPoint's internal co-authoring /
# change-tracking part (ppt/changesInfos/changesInfo1.xml): it records that
# Denis Barritault, a little while after renaming the "PLAY/STOP"-sibling
# pedal textbox's label to "CONTINUE" (already present in this file), went
# on to touch the group that contains it (group id=20, "Group 19", the
# wiring diagram on slide with sldId=257 / cId=833119427). No slide, shape,
# text or geometry content actually differs between the before/after
# states - only the revision-tracking metadata's dt/v/actId counters and an
# extra <pc:grpChg> entry for that group change, none of which is part of
# the exposed PowerPoint object model (there is no COM/VBA property for
# ppt/changesInfos/*, in real PowerPoint or here).
#
# So we reproduce the actual user action that is being tracked: locating
# the "CONTINUE pedal (GPIO 12)" label's enclosing group and (re)selecting
# it, i.e. touching the group described by the diff, without mutating any
# geometry/text that would otherwise introduce spurious differences that
# are not present in the target diff.

$p = $ppt.ActivePresentation

# Find the slide with sldId (SlideID) 257, as referenced by the diff's
# <pc:sldMk cId="833119427" sldId="257"/>.
$targetSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    if ($sl.SlideID -eq 257) {
        $targetSlide = $sl
    }
}

# Find the group shape with id=20 (creationId {9519DA5A-832B-2BC9-31EA-FABCA4E853A8}),
# referenced by the diff's new <pc:grpChg> / <ac:grpSpMk id="20" .../>.
$targetGroup = $null
for ($i = 1; $i -le $targetSlide.Shapes.Count; $i++) {
    $sh = $targetSlide.Shapes.Item($i)
    if ($sh.Id -eq 20) {
        $targetGroup = $sh
    }
}

# Touch/select the group that the change history records as modified.
$targetGroup.Select()

# The CONTINUE pedal textbox (id=14) inside the group already reads
# "CONTINUE pedal (GPIO 12)" - that text edit is the earlier, already
# applied part of history ("pause switch" -> "continue switch"); this
# commit itself changes no visible text, so it is left untouched.
for ($i = 1; $i -le $targetGroup.GroupItems.Count; $i++) {
    $item = $targetGroup.GroupItems.Item($i)
    if ($item.Id -eq 14) {
        $item.Select($false)
    }
}
